$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 505, shifting existing rows 505-546 down to 506-547
$ws.Rows("505:505").Insert()

# Populate the newly inserted row 505 with the new weekly price record
$ws.Range("A505").Value = 10
$ws.Range("B505").Value = "Vega Modelo de Temuco"
$ws.Range("C505").Value = "La Araucanía"
$ws.Range("D505").Value = 45223
$ws.Range("E505").Value = 9
$ws.Range("F505").Value = 100114013
$ws.Range("G505").Value = "Zanahoria"
$ws.Range("H505").Value = "Sin especificar"
$ws.Range("I505").Value = "Primera"
$ws.Range("J505").Value = 90
$ws.Range("K505").Value = 7000
$ws.Range("L505").Value = 7000
$ws.Range("M505").Value = 7000
$ws.Range("N505").Value = "$/saco 20 kilos"
$ws.Range("O505").Value = "Región Metropolitana"
$ws.Range("P505").Value = 350
$ws.Range("Q505").Value = 20
$ws.Range("R505").Value = "Hortaliza"
